$wb = $excel.ActiveWorkbook

# --- Sheet: ManageOrdersPage ---
# Decrease the Order Id values in column A by 3 (new TestNG test data rows added
# earlier in the real table, pushing these order ids down).
$wsOrders = $wb.Worksheets.Item("ManageOrdersPage")
$wsOrders.Range("A2").Value = 145
$wsOrders.Range("A3").Value = 144
$wsOrders.Range("A4").Value = 143
$wsOrders.Range("A5").Value = 142
$wsOrders.Range("A6").Value = 141
$wsOrders.Range("A7").Value = 140
$wsOrders.Range("A8").Value = 139
$wsOrders.Range("A9").Value = 138

# Column A now needs an explicit, best-fit width like the other columns on this sheet
# (Excel auto-fit column A to the "Order Id" header -> stored width 12).
$wsOrders.Columns.Item(1).ColumnWidth = 11.17

# --- Sheet: MobileSliderPage (previously the active/selected tab) ---
# Move the selection off this sheet; it is no longer the tab shown when the
# workbook is reopened.
$wsMobile = $wb.Worksheets.Item("MobileSliderPage")
$wsMobile.Range("B1:B2").Select()

# --- Re-select / re-activate ManageOrdersPage ---
# This becomes the active tab (and its selection moves to A9) when the
# workbook is saved.
$wsOrders.Activate()
$wsOrders.Range("A9").Select()
